# Apply crypto price/volume updates scraped on Sun Mar  3 13:12:23 UTC 2024
# (includes a few coin-ranking swaps where two adjacent rows trade places)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '62.359.38'
$ws.Range("E2").Value = '  +0.53%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.425.10'
$ws.Range("E3").Value = '  +0.21%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  -0.22%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '414.22'
$ws.Range("E5").Value = '  +0.84%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '128.37'
$ws.Range("E6").Value = '  -0.78%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.627'
$ws.Range("E7").Value = '  -2.63%  '

$ws.Range("E8").Value = '  +0.02%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.730'
$ws.Range("E9").Value = '  -1.15%  '

$ws.Range("E10").Value = '  -0.46%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '42.69'
$ws.Range("E11").Value = '  -0.21%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0000220'
$ws.Range("E12").Value = '  +1.26%  '

$ws.Range("E13").Value = '  +0.97%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '3.963.57'
$ws.Range("E14").Value = '  +0.22%  '

$ws.Range("E15").Value = '  -0.31%  '

$ws.Range("E16").Value = '  -3.06%  '

# row 17 -> WrappedEther
$ws.Range("B17").Value = 'WrappedEther'
$ws.Range("C17").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.428.11'
$ws.Range("E17").Value = '  +0.99%  '

# row 18 -> Uniswap
$ws.Range("B18").Value = 'Uniswap'
$ws.Range("C18").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '12.82'
$ws.Range("E18").Value = '  +6.14%  '

$ws.Range("E19").Value = '  -0.93%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '62.268.20'
$ws.Range("E20").Value = '  +0.39%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '477.87'
$ws.Range("E21").Value = '  +7.31%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '92.12'
$ws.Range("E22").Value = '  +0.80%  '

$ws.Range("E23").Value = '  +2.76%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '13.08'
$ws.Range("E24").Value = '  -0.44%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '3.31'
$ws.Range("E25").Value = '  +1.39%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '9.62'
$ws.Range("E26").Value = '  +8.76%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '33.53'
$ws.Range("E27").Value = '  +0.10%  '

$ws.Range("E28").Value = '  +0.27%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.71'
$ws.Range("E29").Value = '  +1.30%  '

# row 30 -> Toncoin
$ws.Range("B30").Value = 'Toncoin'
$ws.Range("C30").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.66'
$ws.Range("E30").Value = '  -3.05%  '

# row 31 -> Cosmos
$ws.Range("B31").Value = 'Cosmos'
$ws.Range("C31").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '11.87'
$ws.Range("E31").Value = '  -0.88%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.167'
$ws.Range("E32").Value = '  -1.55%  '

$ws.Range("E33").Value = '  -3.21%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '41.07'
$ws.Range("E34").Value = '  -4.39%  '

$ws.Range("E35").Value = '  -0.02%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '58.03'
$ws.Range("E36").Value = '  +7.76%  '

$ws.Range("E37").Value = '  -2.71%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.998'
$ws.Range("E38").Value = '  -0.05%  '

$ws.Range("E39").Value = '  +5.04%  '

$ws.Range("E40").Value = '  -0.88%  '

# row 41 -> TheGraph
$ws.Range("B41").Value = 'TheGraph'
$ws.Range("C41").Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.323'
$ws.Range("E41").Value = '  +2.30%  '

# row 42 -> Monero
$ws.Range("B42").Value = 'Monero'
$ws.Range("C42").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '147.42'
$ws.Range("E42").Value = '  +4.46%  '

$ws.Range("E43").Value = '  -1.53%  '

# row 44 -> WEMIXToken
$ws.Range("B44").Value = 'WEMIXToken'
$ws.Range("C44").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.65'
$ws.Range("E44").Value = '  +10.06%  '

# row 45 -> ARBITRUM
$ws.Range("B45").Value = 'ARBITRUM'
$ws.Range("C45").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.07'
$ws.Range("E45").Value = '  +4.87%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '4.28'
$ws.Range("E46").Value = '  +1.21%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.31'
$ws.Range("E47").Value = '  +16.78%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '16.34'
$ws.Range("E48").Value = '  -1.83%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0₃0537'
$ws.Range("E49").Value = '  +26.04%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '22.19'
$ws.Range("E50").Value = '  +0.02%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '113.73'
$ws.Range("E51").Value = '  +8.45%  '
